$d = $word.ActiveDocument

$wdFindContinue = 1
$wdReplaceAll = 2

# --- Edit 1: "...biasing future movements to be more similar to repeated ones..." ---
# "be more similar" -> "become similar"
$d.Content.Find.Execute("be more similar", $true, $false, $false, $false, $false, $true, `
    $wdFindContinue, $false, "become similar", $wdReplaceAll) | Out-Null

# --- Edit 2: "...use-dependent learning in walking from a mechanistic perspective..." ---
# "mechanistic" -> "behavioral"
$d.Content.Find.Execute("mechanistic", $true, $false, $false, $false, $false, $true, `
    $wdFindContinue, $false, "behavioral", $wdReplaceAll) | Out-Null

# --- Edit 3: reorder "whether adaptive or maladaptive," around "gait patterns" ---
# Before: "...contributes to altered, whether adaptive or maladaptive, gait patterns remains unknown..."
# After:  "...contributes to altered gait patterns, whether adaptive or maladaptive,  remains unknown..."
$oldPhrase = ", whether adaptive or maladaptive, gait patterns remains unknown"
$newPhrase = " gait patterns, whether adaptive or maladaptive,  remains unknown"
$d.Content.Find.Execute($oldPhrase, $true, $false, $false, $false, $false, $true, `
    $wdFindContinue, $false, $newPhrase, $wdReplaceAll) | Out-Null

# Word stamps a "_GoBack" bookmark at the site of the most recent edit; place it
# between "gait patterns, " and "whether adaptive or maladaptive, " to match.
$rng = $d.Content
$rng.Find.Execute("whether adaptive or maladaptive,  remains unknown") | Out-Null
$goBackPos = $rng.Start
$goBackRng = $d.Range($goBackPos, $goBackPos)
$d.Bookmarks.Add("_GoBack", $goBackRng) | Out-Null
